$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated CI run counts (more runs accumulated since last snapshot)
$ws.Range("B1").Value = 1106
$ws.Range("B2").Value = 126
$ws.Range("B4").Value = 142
$ws.Range("B5").Value = 37
$ws.Range("B6").Value = 32

# New subtotal formula for the error-breakdown rows (doc/side-effect-first option addition)
$ws.Range("C4").Formula = "=SUM(B4:B8)"

# Leave the cursor where the author left it after editing
[void]$ws.Range("E4").Select()
